$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (separation, average.distance, median.distance, avg.silwidth, average.toother)
$ws.Range("D2").Value = 0.618
$ws.Range("E2").Value = 1.349
$ws.Range("F2").Value = 1.298
$ws.Range("G2").Value = 0.28
$ws.Range("H2").Value = 1.949

# Row 3 updates
$ws.Range("D3").Value = 0.618
$ws.Range("E3").Value = 1.141
$ws.Range("F3").Value = 1.328
$ws.Range("G3").Value = 0.409
$ws.Range("H3").Value = 1.949
